$d = $word.ActiveDocument

# Locate the bibliography paragraph that currently holds the single long run.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("TAYLOR, J. R. - Mec")) {
        $target = $p.Range
        break
    }
}

if ($null -eq $target) {
    throw "Could not locate target paragraph"
}

# Trim the trailing paragraph mark from the range so InsertXML only
# replaces the paragraph's run content (not the paragraph mark itself).
$replaceRange = $d.Range($target.Start, $target.End - 1)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
'<w:p><w:r>' + `
'<w:t xml:space="preserve">TAYLOR, J. R. - Mec&#226;nica Cl&#225;ssica, Bookman, 2015. </w:t><w:br/>' + `
'<w:t xml:space="preserve">THORNTON, S. T. MARION, J. B. &#8211; Din&#226;mica Cl&#225;ssica de Part&#237;culas e Sistemas, tradu&#231;&#227;o da 5&#170; edi&#231;&#227;o norte-americana, CENGAGE Learning, 2016. </w:t><w:br/>' + `
'<w:t xml:space="preserve">F.P. BEER, E.R. JOHNSTON, E. RUSSEL. - Mec&#226;nica vetorial para engenheiros: Est&#225;tica, McGraw Hill. 9a Ed., 2012. </w:t><w:br/>' + `
'<w:t xml:space="preserve">BEER, F.P., JOHNSTON Jr., E.R., CLAUSEN, W. E. - Mec&#226;nica Vetorial para Engenheiros: Din&#226;mica, McGraw-Hill. 7&#170; Ed., 2006. </w:t><w:br/>' + `
'<w:t>GOLDSTEIN, H.; POOLE, C.; SAFKO, J. &#8211; Classical Mechanics, Addison-Wesley Pub. Co. 2013.</w:t><w:br/>' + `
'<w:t>LEMOS, N. A. &#8211; Mec&#226;nica Anal&#237;tica, Livraria da F&#237;sica. 2007.</w:t><w:br/>' + `
'<w:t xml:space="preserve">KOMPANEYETS, A. S. &#8211; Theoretical Physics, Peace Publishers. 2012. </w:t><w:br/>' + `
'<w:t>LANDAU, L. D.; LIFSHITZ, E. M. &#8211; Mechanics, Pergamon Press. 1969</w:t>' + `
'</w:r></w:p>' + `
'</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$replaceRange.InsertXML($xml)

Write-Host "Bibliography paragraph split into" (7 + 1) "lines via manual line breaks."
